$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header cells: Client ID (L1) and Client Secret (M1)
$ws.Range("L1").Value = "Client ID"
$ws.Range("M1").Value = "Client Secret"

# Match the existing header formatting (bold, centered) used by A1:K1
$ws.Range("L1:M1").Font.Bold = $true
$ws.Range("L1:M1").HorizontalAlignment = -4108

# Set the new column widths to match the template
$ws.Columns.Item(12).ColumnWidth = 19.325
$ws.Columns.Item(13).ColumnWidth = 21.166666666666668

# Update the active selection as recorded in the saved file
$ws.Range("J12").Select()
